# Personal Bot - Apresentação
# Fix the typo "Feitas" -> "Freitas" in the author name on the
# title slide's subtitle placeholder ("Leonardo Gomes de Freitas").
#
# The other hunks in the source diff (err="1" spell-check flags being
# cleared, cached datetimeFigureOut placeholder text on slides that
# aren't part of this trimmed deck, and xmlns="" artifacts introduced
# by the original OOXML serializer) are cosmetic round-trip artifacts
# with no visible/semantic effect, and aren't reachable through the
# PowerPoint COM object model exposed here, so they're intentionally
# left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }

    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf("Feitas")

    if ($idx -ge 0) {
        $start = $idx + 1   # TextRange.Characters is 1-indexed
        $sub = $tr.Characters($start, 6)
        $sub.Text = "Freitas"
    }
}
